$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend the 0..13 sequence in B1:O1 with 14 (P1) and 15 (Q1) ---
# Copy the formatting of the last two header cells (N1:O1) into the new P1:Q1
# cells so the new header cells pick up the same bold/border/centered style
# (style index 1) instead of the default style, then set their values.
$ws.Range("N1:O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows (2-25): swap the I/K/M/O values and append the new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column, value 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column, value 2
}
